# Append 12 more rows (206-217) of Korean holiday data, continuing the
# existing A/B series (A = index 204..215, B = fraction value).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(204, 0.8333333333333334),
    @(205, 0.7474747474747474),
    @(206, 0.5777777777777777),
    @(207, 0.8148148148148148),
    @(208, 0.5277777777777778),
    @(209, 0.5396825396825397),
    @(210, 0.7777777777777779),
    @(211, 0.2444444444444444),
    @(212, 0.4444444444444444),
    @(213, 0.2777777777777777),
    @(214, 0.4444444444444444),
    @(215, 0.4444444444444444)
)

$startRow = 206

# Use the existing A205 cell's formatting as the template for the new A cells.
$templateA = $ws.Range("A205")

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $aVal = $data[$i][0]
    $bVal = $data[$i][1]

    $aCell = $ws.Cells.Item($r, 1)
    $bCell = $ws.Cells.Item($r, 2)

    $aCell.Value = $aVal
    $bCell.Value = $bVal

    $templateA.Copy() | Out-Null
    $aCell.PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0
